$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at 25 (between "Ban User" row 24 and the old
#        "Tehnical" header, which was row 26) ---------------------------
$ws.Rows("25:25").Insert()

# Copy the wrap-text / unsolved-red formatting from an existing
# functionality row (row 19) onto the new row 25 so the styles are
# reused instead of duplicated.
$ws.Range("B19:C19").Copy()
$ws.Range("B25:C25").PasteSpecial(-4122)
$ws.Range("A25").Value2 = 24
$ws.Range("B25").Value2 = "Recommend products Random sorted by popularity, newest`nrandom"
$ws.Rows("25:25").RowHeight = 29.4

# --- 2. Insert two more rows before the "Tehnical" header (now row 27)
#        so it ends up at row 29, leaving a 3-row gap (26, 27 spacer, 28
#        blank) above it -------------------------------------------------
$ws.Rows("27:27").Insert()
$ws.Rows("27:27").Insert()

# Row 26 becomes a brand new (not-yet-styled) row; fill it with the new
# functionality + the "solved" green marker, reusing row 21's styling
# for the status cell.
$ws.Range("C21").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("A26").Value2 = 25
$ws.Range("B26").Value2 = "Enhance user login methods, etc"

# Row 27 is a thin blank spacer row below the first table; give C27 the
# white-on-nothing font used for that gap.
$ws.Range("C27").Font.ThemeColor = 2

# --- 3. Append a new "Technical" functionality row at the bottom of the
#        sheet (row 38) --------------------------------------------------
$ws.Range("B35:C35").Copy()
$ws.Range("B38:C38").PasteSpecial(-4122)
$ws.Range("A38").Value2 = 9
$ws.Range("B38").Value2 = "Process data in batches"
$ws.Rows("38:38").RowHeight = 15

$ws.Range("A1").Select()
